$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 1764279.1
$ws.Cells.Item(33, 10).Value = 917
$ws.Cells.Item(33, 12).Value = 917
$ws.Cells.Item(33, 14).Value = -1375
$ws.Cells.Item(69, 8).Value = 6137.9707
$ws.Cells.Item(69, 9).Value = 3991.6667
$ws.Cells.Item(69, 10).Value = 6597.893
$ws.Cells.Item(69, 11).Value = 11975.0001
$ws.Cells.Item(69, 12).Value = 19793.679
$ws.Cells.Item(69, 13).Value = -11101.0001
$ws.Cells.Item(69, 14).Value = -21541.679
$ws.Cells.Item(72, 8).Value = 6137.9707
$ws.Cells.Item(72, 9).Value = 3991.6667
$ws.Cells.Item(72, 10).Value = 6597.893
$ws.Cells.Item(72, 11).Value = 35925.0003
$ws.Cells.Item(72, 12).Value = 59381.037
$ws.Cells.Item(72, 13).Value = -31557.0003
$ws.Cells.Item(72, 14).Value = -68117.037
$ws.Cells.Item(75, 8).Value = 45500
$ws.Cells.Item(75, 9).Value = 40000
$ws.Cells.Item(75, 10).Value = 51000
$ws.Cells.Item(75, 11).Value = 40000
$ws.Cells.Item(75, 12).Value = 51000
$ws.Cells.Item(75, 13).Value = -39064
$ws.Cells.Item(75, 14).Value = -52872
$ws.Cells.Item(76, 8).Value = 5553.095
$ws.Cells.Item(76, 9).Value = 4556.5386
$ws.Cells.Item(76, 11).Value = 4556.5386
$ws.Cells.Item(76, 13).Value = -4241.5386
$ws.Cells.Item(78, 8).Value = 45500
$ws.Cells.Item(78, 9).Value = 40000
$ws.Cells.Item(78, 10).Value = 51000
$ws.Cells.Item(78, 11).Value = 120000
$ws.Cells.Item(78, 12).Value = 153000
$ws.Cells.Item(78, 13).Value = -115320
$ws.Cells.Item(78, 14).Value = -162360
$ws.Cells.Item(79, 8).Value = 5553.095
$ws.Cells.Item(79, 9).Value = 4556.5386
$ws.Cells.Item(79, 11).Value = 4556.5386
$ws.Cells.Item(79, 13).Value = -3464.5386
$ws.Cells.Item(100, 8).Value = 859.7222
$ws.Cells.Item(100, 9).Value = 876.63635
$ws.Cells.Item(100, 10).Value = 833.1429000000001
$ws.Cells.Item(100, 11).Value = 876.63635
$ws.Cells.Item(100, 12).Value = 833.1429000000001
$ws.Cells.Item(100, 13).Value = -335.63635
$ws.Cells.Item(100, 14).Value = -1915.1429
$ws.Cells.Item(116, 8).Value = 5510.6895
$ws.Cells.Item(116, 9).Value = 3491.8333
$ws.Cells.Item(116, 10).Value = 6935.7646
$ws.Cells.Item(116, 11).Value = 3491.8333
$ws.Cells.Item(116, 12).Value = 6935.7646
$ws.Cells.Item(116, 13).Value = -49.83329999999978
$ws.Cells.Item(116, 14).Value = -13819.7646
$ws.Cells.Item(118, 8).Value = 90909630
$ws.Cells.Item(118, 9).Value = 142857540
$ws.Cells.Item(118, 10).Value = 787
$ws.Cells.Item(118, 11).Value = 428572620
$ws.Cells.Item(118, 12).Value = 2361
$ws.Cells.Item(118, 13).Value = -428570963
$ws.Cells.Item(118, 14).Value = -5675
$ws.Cells.Item(127, 8).Value = 1162.8334
$ws.Cells.Item(127, 9).Value = 578.3333
$ws.Cells.Item(127, 10).Value = 2331.8333
$ws.Cells.Item(127, 11).Value = 1734.9999
$ws.Cells.Item(127, 12).Value = 6995.499899999999
$ws.Cells.Item(127, 13).Value = 3225.0001
$ws.Cells.Item(127, 14).Value = -16915.4999
$ws.Cells.Item(132, 8).Value = 5193.625
$ws.Cells.Item(132, 9).Value = 6055.8423
$ws.Cells.Item(132, 10).Value = 1917.2
$ws.Cells.Item(132, 11).Value = 18167.5269
$ws.Cells.Item(132, 12).Value = 5751.6
$ws.Cells.Item(132, 13).Value = -15637.5269
$ws.Cells.Item(132, 14).Value = -10811.6
$ws.Cells.Item(137, 8).Value = 59635.87
$ws.Cells.Item(137, 9).Value = 106296.94
$ws.Cells.Item(137, 11).Value = 318890.82
$ws.Cells.Item(137, 13).Value = -316340.82
$ws.Cells.Item(138, 8).Value = 3289.0679
$ws.Cells.Item(138, 10).Value = 3555.7727
$ws.Cells.Item(138, 12).Value = 10667.3181
$ws.Cells.Item(138, 14).Value = -20947.3181

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 943309.9399999999
$ws.Cells.Item(2, 9).Value = 1047947.25
$ws.Cells.Item(2, 11).Value = 1047947.25
$ws.Cells.Item(2, 13).Value = -1047834.25
$ws.Cells.Item(45, 8).Value = 4331574
$ws.Cells.Item(45, 9).Value = 6494850.5
$ws.Cells.Item(45, 10).Value = 5021.5454
$ws.Cells.Item(45, 11).Value = 6494850.5
$ws.Cells.Item(45, 12).Value = 5021.5454
$ws.Cells.Item(45, 13).Value = -6494473.5
$ws.Cells.Item(45, 14).Value = -5775.5454
$ws.Cells.Item(61, 8).Value = 4238.1665
$ws.Cells.Item(61, 9).Value = 4087.1333
$ws.Cells.Item(61, 11).Value = 4087.1333
$ws.Cells.Item(61, 13).Value = -3875.1333
$ws.Cells.Item(74, 8).Value = 141739
$ws.Cells.Item(74, 9).Value = 130029.43
$ws.Cells.Item(74, 10).Value = 155400.17
$ws.Cells.Item(74, 11).Value = 130029.43
$ws.Cells.Item(74, 12).Value = 155400.17
$ws.Cells.Item(74, 13).Value = -129155.43
$ws.Cells.Item(74, 14).Value = -157148.17
$ws.Cells.Item(77, 8).Value = 141739
$ws.Cells.Item(77, 9).Value = 130029.43
$ws.Cells.Item(77, 10).Value = 155400.17
$ws.Cells.Item(77, 11).Value = 650147.1499999999
$ws.Cells.Item(77, 12).Value = 777000.8500000001
$ws.Cells.Item(77, 13).Value = -645779.1499999999
$ws.Cells.Item(77, 14).Value = -785736.8500000001
$ws.Cells.Item(97, 8).Value = 1044534.06
$ws.Cells.Item(97, 9).Value = 1285413.9
$ws.Cells.Item(97, 10).Value = 721.3333
$ws.Cells.Item(97, 11).Value = 1285413.9
$ws.Cells.Item(97, 12).Value = 721.3333
$ws.Cells.Item(97, 13).Value = -1284917.9
$ws.Cells.Item(97, 14).Value = -1713.3333
$ws.Cells.Item(116, 8).Value = 943309.9399999999
$ws.Cells.Item(116, 9).Value = 1047947.25
$ws.Cells.Item(116, 11).Value = 1047947.25
$ws.Cells.Item(116, 13).Value = -1045653.25
$ws.Cells.Item(122, 8).Value = 2828973.8
$ws.Cells.Item(122, 9).Value = 3291174
$ws.Cells.Item(122, 10).Value = 2089452.9
$ws.Cells.Item(122, 11).Value = 9873522
$ws.Cells.Item(122, 12).Value = 6268358.699999999
$ws.Cells.Item(122, 13).Value = -9871072
$ws.Cells.Item(122, 14).Value = -6273258.699999999
$ws.Cells.Item(132, 8).Value = 2579.0334
$ws.Cells.Item(132, 9).Value = 1635.2106
$ws.Cells.Item(132, 10).Value = 4209.273
$ws.Cells.Item(132, 11).Value = 4905.6318
$ws.Cells.Item(132, 12).Value = 12627.819
$ws.Cells.Item(132, 13).Value = -2375.6318
$ws.Cells.Item(132, 14).Value = -17687.819
$ws.Cells.Item(136, 8).Value = 4238.1665
$ws.Cells.Item(136, 9).Value = 4087.1333
$ws.Cells.Item(136, 11).Value = 12261.3999
$ws.Cells.Item(136, 13).Value = -9711.3999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 943309.9399999999
$ws.Cells.Item(3, 9).Value = 1047947.25
$ws.Cells.Item(3, 11).Value = 1047947.25
$ws.Cells.Item(3, 13).Value = -1047833.25
$ws.Cells.Item(80, 8).Value = 629
$ws.Cells.Item(80, 10).Value = 558.2
$ws.Cells.Item(80, 12).Value = 558.2
$ws.Cells.Item(80, 14).Value = -2554.2
$ws.Cells.Item(83, 8).Value = 629
$ws.Cells.Item(83, 10).Value = 558.2
$ws.Cells.Item(83, 12).Value = 2791
$ws.Cells.Item(83, 14).Value = -12775
$ws.Cells.Item(94, 8).Value = 4004895.8
$ws.Cells.Item(94, 9).Value = 4547200
$ws.Cells.Item(94, 10).Value = 27998.666
$ws.Cells.Item(94, 11).Value = 4547200
$ws.Cells.Item(94, 12).Value = 27998.666
$ws.Cells.Item(94, 13).Value = -4546749
$ws.Cells.Item(94, 14).Value = -28900.666
$ws.Cells.Item(99, 8).Value = 6213161
$ws.Cells.Item(99, 9).Value = 7520116.5
$ws.Cells.Item(99, 10).Value = 5124.25
$ws.Cells.Item(99, 11).Value = 7520116.5
$ws.Cells.Item(99, 12).Value = 5124.25
$ws.Cells.Item(99, 13).Value = -7518618.5
$ws.Cells.Item(99, 14).Value = -8120.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 337.84616
$ws.Cells.Item(22, 9).Value = 325.14285
$ws.Cells.Item(22, 10).Value = 352.66666
$ws.Cells.Item(22, 11).Value = 325.14285
$ws.Cells.Item(22, 12).Value = 352.66666
$ws.Cells.Item(22, 13).Value = 24.85714999999999
$ws.Cells.Item(22, 14).Value = -1052.66666
$ws.Cells.Item(31, 8).Value = 21898.04
$ws.Cells.Item(31, 9).Value = 3229.8823
$ws.Cells.Item(31, 11).Value = 3229.8823
$ws.Cells.Item(31, 13).Value = -2934.8823
$ws.Cells.Item(34, 8).Value = 21898.04
$ws.Cells.Item(34, 9).Value = 3229.8823
$ws.Cells.Item(34, 11).Value = 3229.8823
$ws.Cells.Item(34, 13).Value = -3027.8823
$ws.Cells.Item(99, 8).Value = 3427.625
$ws.Cells.Item(99, 9).Value = 3559.9167
$ws.Cells.Item(99, 10).Value = 3030.75
$ws.Cells.Item(99, 11).Value = 3559.9167
$ws.Cells.Item(99, 12).Value = 3030.75
$ws.Cells.Item(99, 13).Value = -2061.9167
$ws.Cells.Item(99, 14).Value = -6026.75
$ws.Cells.Item(126, 8).Value = 3427.625
$ws.Cells.Item(126, 9).Value = 3559.9167
$ws.Cells.Item(126, 10).Value = 3030.75
$ws.Cells.Item(126, 11).Value = 10679.7501
$ws.Cells.Item(126, 12).Value = 9092.25
$ws.Cells.Item(126, 13).Value = -8209.750100000001
$ws.Cells.Item(126, 14).Value = -14032.25
$ws.Cells.Item(130, 8).Value = 62920
$ws.Cells.Item(130, 10).Value = 62920
$ws.Cells.Item(130, 12).Value = 62920
$ws.Cells.Item(130, 14).Value = -72960
$ws.Cells.Item(133, 8).Value = 69999
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 69999
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 12).Value = 69999
$ws.Cells.Item(133, 13).ClearContents()
$ws.Cells.Item(133, 14).Value = -75059
$ws.Cells.Item(134, 8).Value = 3676.4443
$ws.Cells.Item(134, 9).Value = 2322
$ws.Cells.Item(134, 11).Value = 6966
$ws.Cells.Item(134, 13).Value = -4431
$ws.Cells.Item(137, 8).Value = 124633
$ws.Cells.Item(137, 10).Value = 124633
$ws.Cells.Item(137, 12).Value = 124633
$ws.Cells.Item(137, 14).Value = -134833

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 56850.945
$ws.Cells.Item(5, 10).Value = 126919.875
$ws.Cells.Item(5, 12).Value = 380759.625
$ws.Cells.Item(5, 14).Value = -380983.625
$ws.Cells.Item(34, 8).Value = 1819.8
$ws.Cells.Item(34, 10).Value = 1819.8
$ws.Cells.Item(34, 12).Value = 5459.4
$ws.Cells.Item(34, 14).Value = -5627.4
$ws.Cells.Item(37, 8).Value = 47990
$ws.Cells.Item(37, 10).Value = 47990
$ws.Cells.Item(37, 12).Value = 143970
$ws.Cells.Item(37, 14).Value = -144194
$ws.Cells.Item(52, 8).Value = 999
$ws.Cells.Item(52, 10).Value = 999
$ws.Cells.Item(52, 12).Value = 2997
$ws.Cells.Item(52, 14).Value = -3529
$ws.Cells.Item(113, 8).Value = 3478.8262
$ws.Cells.Item(113, 10).Value = 1697.5
$ws.Cells.Item(113, 12).Value = 5092.5
$ws.Cells.Item(113, 14).Value = -9432.5
$ws.Cells.Item(126, 8).Value = 2971.3333
$ws.Cells.Item(126, 9).Value = 2442.6667
$ws.Cells.Item(126, 10).Value = 3500
$ws.Cells.Item(126, 11).Value = 7328.000100000001
$ws.Cells.Item(126, 12).Value = 10500
$ws.Cells.Item(126, 13).Value = -2388.000100000001
$ws.Cells.Item(126, 14).Value = -20380
$ws.Cells.Item(132, 8).Value = 2110.9092
$ws.Cells.Item(132, 10).Value = 2387.2727
$ws.Cells.Item(132, 12).Value = 21485.4543
$ws.Cells.Item(132, 14).Value = -26545.4543
$ws.Cells.Item(135, 8).Value = 56850.945
$ws.Cells.Item(135, 10).Value = 126919.875
$ws.Cells.Item(135, 12).Value = 1142278.875
$ws.Cells.Item(135, 14).Value = -1147348.875

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(12, 8).Value = 27555.445
$ws.Cells.Item(12, 10).Value = 48999.5
$ws.Cells.Item(12, 12).Value = 48999.5
$ws.Cells.Item(12, 14).Value = -49279.5
$ws.Cells.Item(21, 8).Value = 3900
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 12).Value = 0
$ws.Cells.Item(21, 14).ClearContents()
$ws.Cells.Item(30, 8).Value = 3900
$ws.Cells.Item(30, 10).Value = 0
$ws.Cells.Item(30, 12).Value = 0
$ws.Cells.Item(30, 14).ClearContents()
$ws.Cells.Item(80, 8).Value = 13457739
$ws.Cells.Item(80, 9).Value = 20177588
$ws.Cells.Item(80, 11).Value = 20177588
$ws.Cells.Item(80, 13).Value = -20176590
$ws.Cells.Item(83, 8).Value = 13457739
$ws.Cells.Item(83, 9).Value = 20177588
$ws.Cells.Item(83, 11).Value = 100887940
$ws.Cells.Item(83, 13).Value = -100882948
$ws.Cells.Item(97, 8).Value = 23809524
$ws.Cells.Item(97, 9).Value = 23809524
$ws.Cells.Item(97, 11).Value = 23809524
$ws.Cells.Item(97, 13).Value = -23809028
$ws.Cells.Item(99, 8).Value = 9323.75
$ws.Cells.Item(99, 9).Value = 9323.75
$ws.Cells.Item(99, 11).Value = 9323.75
$ws.Cells.Item(99, 13).Value = -7077.75
$ws.Cells.Item(113, 8).Value = 4067132.2
$ws.Cells.Item(113, 9).Value = 5557047
$ws.Cells.Item(113, 10).Value = 3728.3635
$ws.Cells.Item(113, 11).Value = 5557047
$ws.Cells.Item(113, 12).Value = 3728.3635
$ws.Cells.Item(113, 13).Value = -5554877
$ws.Cells.Item(113, 14).Value = -8068.363499999999
$ws.Cells.Item(122, 8).Value = 345513.7
$ws.Cells.Item(122, 9).Value = 426088.38
$ws.Cells.Item(122, 10).Value = 7100
$ws.Cells.Item(122, 11).Value = 1278265.14
$ws.Cells.Item(122, 12).Value = 21300
$ws.Cells.Item(122, 13).Value = -1275815.14
$ws.Cells.Item(122, 14).Value = -26200
$ws.Cells.Item(132, 8).Value = 3208.4707
$ws.Cells.Item(132, 9).Value = 3112.1738
$ws.Cells.Item(132, 10).Value = 3409.818
$ws.Cells.Item(132, 11).Value = 9336.5214
$ws.Cells.Item(132, 12).Value = 10229.454
$ws.Cells.Item(132, 13).Value = -6806.5214
$ws.Cells.Item(132, 14).Value = -15289.454
$ws.Cells.Item(136, 8).Value = 11770.656
$ws.Cells.Item(136, 10).Value = 11770.656
$ws.Cells.Item(136, 12).Value = 35311.968
$ws.Cells.Item(136, 14).Value = -40411.968

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 43235.906
$ws.Cells.Item(22, 9).Value = 222527
$ws.Cells.Item(22, 10).Value = 1049.7646
$ws.Cells.Item(22, 11).Value = 222527
$ws.Cells.Item(22, 12).Value = 1049.7646
$ws.Cells.Item(22, 13).Value = -222232
$ws.Cells.Item(22, 14).Value = -1639.7646
$ws.Cells.Item(27, 8).Value = 43235.906
$ws.Cells.Item(27, 9).Value = 222527
$ws.Cells.Item(27, 10).Value = 1049.7646
$ws.Cells.Item(27, 11).Value = 222527
$ws.Cells.Item(27, 12).Value = 1049.7646
$ws.Cells.Item(27, 13).Value = -222420
$ws.Cells.Item(27, 14).Value = -1263.7646
$ws.Cells.Item(43, 8).Value = 0
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 11).Value = 0
$ws.Cells.Item(43, 13).ClearContents()
$ws.Cells.Item(46, 8).Value = 2554.9487
$ws.Cells.Item(46, 9).Value = 1659.2
$ws.Cells.Item(46, 10).Value = 2686.6765
$ws.Cells.Item(46, 11).Value = 1659.2
$ws.Cells.Item(46, 12).Value = 2686.6765
$ws.Cells.Item(46, 13).Value = -1471.2
$ws.Cells.Item(46, 14).Value = -3062.6765
$ws.Cells.Item(55, 8).Value = 1397
$ws.Cells.Item(55, 9).Value = 1043.7
$ws.Cells.Item(55, 10).Value = 1691.4166
$ws.Cells.Item(55, 11).Value = 1043.7
$ws.Cells.Item(55, 12).Value = 1691.4166
$ws.Cells.Item(55, 13).Value = -870.7
$ws.Cells.Item(55, 14).Value = -2037.4166
$ws.Cells.Item(61, 8).Value = 4116091
$ws.Cells.Item(61, 9).Value = 4630430.5
$ws.Cells.Item(61, 11).Value = 4630430.5
$ws.Cells.Item(61, 13).Value = -4630228.5
$ws.Cells.Item(68, 8).Value = 4000
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 14).ClearContents()
$ws.Cells.Item(71, 8).Value = 4000
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 14).ClearContents()
$ws.Cells.Item(93, 8).Value = 17556118
$ws.Cells.Item(93, 9).Value = 25642560
$ws.Cells.Item(93, 10).Value = 35498.668
$ws.Cells.Item(93, 11).Value = 25642560
$ws.Cells.Item(93, 12).Value = 35498.668
$ws.Cells.Item(93, 13).Value = -25641312
$ws.Cells.Item(93, 14).Value = -37994.668
$ws.Cells.Item(113, 8).Value = 4116091
$ws.Cells.Item(113, 9).Value = 4630430.5
$ws.Cells.Item(113, 11).Value = 4630430.5
$ws.Cells.Item(113, 13).Value = -4628260.5
$ws.Cells.Item(122, 8).Value = 5509.905
$ws.Cells.Item(122, 9).Value = 3738.818
$ws.Cells.Item(122, 10).Value = 7458.1
$ws.Cells.Item(122, 11).Value = 11216.454
$ws.Cells.Item(122, 12).Value = 22374.3
$ws.Cells.Item(122, 13).Value = -8766.454000000002
$ws.Cells.Item(122, 14).Value = -27274.3

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(45, 8).Value = 22643.75
$ws.Cells.Item(45, 9).Value = 9955
$ws.Cells.Item(45, 10).Value = 26873.334
$ws.Cells.Item(45, 11).Value = 9955
$ws.Cells.Item(45, 12).Value = 26873.334
$ws.Cells.Item(45, 13).Value = -9464
$ws.Cells.Item(45, 14).Value = -27855.334
$ws.Cells.Item(58, 8).Value = 6124.75
$ws.Cells.Item(58, 9).Value = 6124.75
$ws.Cells.Item(58, 11).Value = 6124.75
$ws.Cells.Item(58, 13).Value = -5816.75
$ws.Cells.Item(62, 8).Value = 7595.7393
$ws.Cells.Item(62, 9).Value = 3751
$ws.Cells.Item(62, 10).Value = 7863.9766
$ws.Cells.Item(62, 11).Value = 3751
$ws.Cells.Item(62, 12).Value = 7863.9766
$ws.Cells.Item(62, 13).Value = -3127
$ws.Cells.Item(62, 14).Value = -9111.9766
$ws.Cells.Item(65, 8).Value = 7595.7393
$ws.Cells.Item(65, 9).Value = 3751
$ws.Cells.Item(65, 10).Value = 7863.9766
$ws.Cells.Item(65, 11).Value = 18755
$ws.Cells.Item(65, 12).Value = 39319.883
$ws.Cells.Item(65, 13).Value = -15635
$ws.Cells.Item(65, 14).Value = -45559.883
$ws.Cells.Item(70, 8).Value = 55050
$ws.Cells.Item(70, 9).Value = 20000
$ws.Cells.Item(70, 11).Value = 20000
$ws.Cells.Item(70, 13).Value = -19685
$ws.Cells.Item(73, 8).Value = 55050
$ws.Cells.Item(73, 9).Value = 20000
$ws.Cells.Item(73, 11).Value = 20000
$ws.Cells.Item(73, 13).Value = -18908
$ws.Cells.Item(96, 8).Value = 2689.6
$ws.Cells.Item(96, 9).Value = 2689.6
$ws.Cells.Item(96, 11).Value = 2689.6
$ws.Cells.Item(96, 13).Value = -1316.6
$ws.Cells.Item(107, 8).Value = 50001100
$ws.Cells.Item(107, 9).Value = 100000850
$ws.Cells.Item(107, 10).Value = 1350.8
$ws.Cells.Item(107, 11).Value = 300002550
$ws.Cells.Item(107, 12).Value = 4052.4
$ws.Cells.Item(107, 13).Value = -300000630
$ws.Cells.Item(107, 14).Value = -7892.4
$ws.Cells.Item(132, 8).Value = 15091082
$ws.Cells.Item(132, 9).Value = 18521866
$ws.Cells.Item(132, 10).Value = 840131.25
$ws.Cells.Item(132, 11).Value = 55565598
$ws.Cells.Item(132, 12).Value = 2520393.75
$ws.Cells.Item(132, 13).Value = -55563068
$ws.Cells.Item(132, 14).Value = -2525453.75
$ws.Cells.Item(135, 8).Value = 51905
$ws.Cells.Item(135, 10).Value = 51905
$ws.Cells.Item(135, 12).Value = 51905
$ws.Cells.Item(135, 14).Value = -62045
$ws.Cells.Item(136, 8).Value = 1144.1628
$ws.Cells.Item(136, 9).Value = 1014.1389
$ws.Cells.Item(136, 10).Value = 1812.8572
$ws.Cells.Item(136, 11).Value = 3042.4167
$ws.Cells.Item(136, 12).Value = 5438.571599999999
$ws.Cells.Item(136, 13).Value = -492.4167000000002
$ws.Cells.Item(136, 14).Value = -10538.5716
